$d = $word.ActiveDocument

function Escape-Xml($s) {
    $s = $s -replace '&', '&amp;'
    $s = $s -replace '<', '&lt;'
    $s = $s -replace '>', '&gt;'
    return $s
}

# Replaces the text of a single run located at character offset $start (length
# of $old) with $new, optionally carrying run properties ($rPrXml, e.g.
# "<w:rPr><w:b/></w:rPr>"). Uses Range.InsertXML (wrapped in a minimal
# w:document/w:body/w:p) instead of Range.Text / Find-Replace so that sibling
# (empty) runs in the same paragraph are left completely untouched.
function Set-RunTextAt($d, $start, $old, $new, $rPrXml = "") {
    $len = $old.Length
    $rng = $d.Range($start, $start + $len)
    if ($rng.Text -ne $old) {
        throw "Range mismatch: expected [$old] but found [$($rng.Text)]"
    }
    $safeNew = Escape-Xml $new
    $xml = '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r>' + $rPrXml + '<w:t>' + $safeNew + '</w:t></w:r></w:p></w:body></w:document>'
    $rng.InsertXML($xml)
}

# Finds the (searchFrom-th, 0-based) occurrence of $old in the document text
# and replaces it using Set-RunTextAt.
function Replace-Nth($d, $old, $new, $occurrenceIndex = 0, $rPrXml = "") {
    $full = $d.Content.Text
    $searchFrom = 0
    $idx = -1
    for ($i = 0; $i -le $occurrenceIndex; $i++) {
        $idx = $full.IndexOf($old, $searchFrom)
        if ($idx -lt 0) {
            throw "Text not found (occurrence $i): $old"
        }
        $searchFrom = $idx + 1
    }
    Set-RunTextAt $d $idx $old $new $rPrXml
}

# --- Heading1 title ---
Replace-Nth $d "Play Bubble Craze for Free: A Unique Slot Game by IGT" "Play Bubble Craze Free and Experience Unique Gameplay" 0

# --- "What we like" bullets ---
Replace-Nth $d "Unique approach to slot machine gameplay" "Unique gameplay mechanics and rules" 0
Replace-Nth $d "Colorful and well-designed graphics" "Excellent graphics and design" 0
Replace-Nth $d "Generous free spins bonus feature" "Free spins bonus feature with guaranteed wins" 0
Replace-Nth $d "Developed by trusted provider IGT" "Developed by a leading provider in the industry" 0

# --- "What we don't like" bullets ---
Replace-Nth $d "Limited bonus feature with no re-triggering" "Bonus round can only be triggered once" 0
Replace-Nth $d "Payouts can be lower than traditional slot machines" "Limited number of bonus features" 0

# --- Bold title repeated near the end ---
Replace-Nth $d "Play Bubble Craze for Free: A Unique Slot Game by IGT" "Play Bubble Craze Free and Experience Unique Gameplay" 0 "<w:rPr><w:b/></w:rPr>"

# --- Italic summary sentence ---
Replace-Nth $d "Try Bubble Craze for free and experience a unique approach to slot machine gameplay with excellent graphics and bonus features by trusted developer IGT." "Read our review of Bubble Craze and play for free to enjoy its unique gameplay mechanics." 0 "<w:rPr><w:i/></w:rPr>"
